$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date column to be treated as plain text, matching the rest of the column
$ws.Range("A8:A9").NumberFormat = "@"

# Row 8: 2026-03-19, Limited, 12800, 2, 1, 1, Auto-generated from bookings
$ws.Range("A8").Value = "2026-03-19"
$ws.Range("B8").Value = "Limited"
$ws.Range("C8").Value = 12800
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "Auto-generated from bookings"

# Row 9: 2026-03-20, Booked, 12800, 2, 2, 0, Auto-generated from bookings
$ws.Range("A9").Value = "2026-03-20"
$ws.Range("B9").Value = "Booked"
$ws.Range("C9").Value = 12800
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = "Auto-generated from bookings"
